$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Objetivos -> now shows the professor name (content that used to live in row 13)
$ws.Range("B10").Value = "5840897 - Clodoaldo Saron"
$ws.Range("C10").Value = "5840897 - Clodoaldo Saron"

# Row 13: used to be the (unlabeled) professor-name row; now becomes "Programa resumido:" / "Semestral"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# Row 14: used to be "Programa resumido:" / long syllabus text; now becomes "Short syllabus:" (no B/C)
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = $null
$ws.Range("C14").Value = $null

# Row 15: used to be "Short syllabus:"; now becomes "Programa:" / "01/01/2012"
# (force text via NumberFormat "@" so the date-like string isn't coerced into
# a date serial, then restore the original General-format style so B15/C15
# end up on the same style indices as the rest of column B/C)
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "01/01/2012"
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Rows.Item(15).RowHeight = 120

# Row 16: used to be "Programa:" / long program text; now becomes "Syllabus:" (no B/C)
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = $null
$ws.Range("C16").Value = $null

# Row 17: used to be "Syllabus:" with a 120pt custom height; now becomes "Avaliação:" with default height
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).EntireRow.AutoFit()

# Row 18: used to be "Avaliação:" with default height; now becomes "Método:" / professor name, 60pt height
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5840897 - Clodoaldo Saron"
$ws.Range("C18").Value = "5840897 - Clodoaldo Saron"
$ws.Rows.Item(18).RowHeight = 60

# Row 19: used to be "Método:" / class-method text; now becomes "Critério:" (B/C text unchanged)
$ws.Range("A19").Value = "Critério:"

# Row 20: used to be "Critério:" / grading-average text; now becomes "Norma de recuperação:" (B/C unchanged)
$ws.Range("A20").Value = "Norma de recuperação:"

# Row 21: used to be "Norma de recuperação:" / makeup-exam text; now becomes "Bibliografia:" (B/C unchanged),
# and the row grows from a 60pt to a 120pt custom height
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120

# Row 22: used to be "Bibliografia:" / bibliography text (120pt); now becomes "Requisitos:" with default height, no B/C
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Value = $null
$ws.Range("C22").Value = $null
$ws.Rows.Item(22).EntireRow.AutoFit()

# Row 23: used to be "Requisitos:" with no B/C; now becomes the requirement text in B/C with no A, 30pt height
$ws.Range("A23").Value = $null
$ws.Range("B23").Value = "LOM3212 -  Fenômenos de Transporte A  (Requisito)`n"
$ws.Range("C23").Value = "LOM3212 -  Fenômenos de Transporte A  (Requisito)`n"
$ws.Rows.Item(23).RowHeight = 30

# Row 24 (old requirement-text row) is now redundant - remove it entirely
$ws.Rows.Item(24).Delete()
